$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift old column G (date_available, with its values+style) into the new column H
$ws.Range("G1:G4").Copy($ws.Range("H1:H4"))

# Fill in the new column G with utilities info (who pays utilities)
$ws.Range("G1").Value = "utilities"
$ws.Range("G2").Value = "landlord"
$ws.Range("G3").Value = "tenant"
$ws.Range("G4").Value = "landlord"

# Give column G the same formatting (centered, style index 1) as the other
# plain-text columns, by copying the format from column F.
$ws.Range("F1:F4").Copy()
$ws.Range("G1:G4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Column widths: G goes back to the standard width (like F); H takes on the
# bestFit width the old "date_available" column used to have.
$ws.Range("G1:G1048576").ColumnWidth = 8
$ws.Range("H1:H1048576").ColumnWidth = 9.6

$ws.Range("H2").Select()
